$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 42
$ws.Range("F8").Value = 909
$ws.Range("F9").Value = 1612
$ws.Range("F10").Value = 6072
$ws.Range("F14").Value = 5922
$ws.Range("F15").Value = 115
$ws.Range("F19").Value = 1647
$ws.Range("F22").Value = 142
$ws.Range("F23").Value = 1376
$ws.Range("F25").Value = 241
$ws.Range("F28").Value = 28
$ws.Range("F30").Value = 3860
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 163
$ws.Range("F8").Value = 377
$ws.Range("F13").Value = 1
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9496
$ws.Range("F4").Value = 619
$ws.Range("F5").Value = 183
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9496
$ws.Range("F4").Value = 619
$ws.Range("F9").Value = 42
$ws.Range("F12").Value = 909
$ws.Range("F13").Value = 183
$ws.Range("F14").Value = 1612
$ws.Range("F15").Value = 6072
$ws.Range("F23").Value = 5922
$ws.Range("F24").Value = 115
$ws.Range("F28").Value = 1647
$ws.Range("F31").Value = 142
$ws.Range("F32").Value = 1376
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 241
$ws.Range("F45").Value = 3860
